# Updated Date and Email checking script to handle emails with special characters.
# This inserts a new "Special Characters Name" row into the InValidEmail sheet
# (pushing the existing "Incorrect Email" row down from row 2 to row 3), and
# updates the InValidEmail defined name to cover the new range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InValidEmail")

# Insert a new blank row above the existing data row (old row 2 -> row 3)
$ws.Rows.Item(2).Insert()

# Make sure the new row's cells are stored as text (matching the rest of the sheet,
# which stores every value - including numeric-looking ones - as text)
$ws.Range("A2:G2").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "131"
$ws.Cells.Item(2, 2).Value = "Special Characters Name"
$ws.Cells.Item(2, 3).Value = "34"
$ws.Cells.Item(2, 4).Value = "special`$`$name@example.com"
$ws.Cells.Item(2, 5).Value = "2021-04-17"
$ws.Cells.Item(2, 6).Value = "2023-02-07"
$ws.Cells.Item(2, 7).Value = "200.5"

# Restore default/normal style on the row that got shifted down
$ws.Range("A3:G3").Style = "Normal"

# Update the InValidEmail defined name range to include the new row
$n = $wb.Names.Item("InValidEmail")
$n.RefersTo = "='InValidEmail'!`$A`$1:`$G`$3"
